# Add two new columns, I (I0) and J (IF), to the existing data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Clone the header style (bold, centered, bordered) from H1 onto the two
# new header cells so the new columns look like the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial($xlPasteFormats)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-11: column I = "I0", column J = "IF".
$iValues = @(9, 7, 8, 9, 9, 3, 2, 6, 7, 8)
$jValues = @(9, 8, 8, 9, 9, 5, 6, 8, 7, 8)

for ($r = 0; $r -lt $iValues.Count; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
